$wb = $excel.ActiveWorkbook

# --- Sheet "pro" (sheet1): update column B values ---
$wsPro = $wb.Worksheets.Item("pro")
$wsPro.Range("B2").Value = 2858042.4441151237
$wsPro.Range("B3").Value = 3234474.5079740775
$wsPro.Range("B4").Value = 3565093.488362123
$wsPro.Range("B5").Value = 3529880.433007928
$wsPro.Range("B6").Value = 3419058.00365659
$wsPro.Range("B7").Value = 3483015.9777827743
$wsPro.Range("B8").Value = 3447876.9166528936
$wsPro.Range("B9").Value = 3260423.887855372
$wsPro.Range("B10").Value = 3437118.5208715466
$wsPro.Range("B11").Value = 3592483.8895612396
$wsPro.Range("B12").Value = 3664033.855631725
$wsPro.Range("B13").Value = 3770294.4123181165
$wsPro.Range("B14").Value = 4019364.579925042
$wsPro.Range("B15").Value = 4142978.93488257
$wsPro.Range("B16").Value = 4498682.490473724
$wsPro.Range("B17").Value = 4381460.758833277
$wsPro.Range("B18").Value = 5131718.955442092
$wsPro.Range("B19").Value = 5378021.941795561
$wsPro.Range("B20").Value = 5722553.458223218
$wsPro.Range("B21").Value = 6097117.0
$wsPro.Range("B22").Value = 6710446.0
$wsPro.Range("B23").Value = 6780085.11732773
$wsPro.Range("B24").Value = 7567508.219889635
$wsPro.Range("B25").Value = 7885018.894669292
$wsPro.Range("B26").Value = 7712467.935935166

# --- Sheet "ind" (sheet2): update column B values ---
$wsInd = $wb.Worksheets.Item("ind")
$wsInd.Range("B2").Value = 936105.7598218526
$wsInd.Range("B3").Value = 928314.3757971823
$wsInd.Range("B4").Value = 932669.0630365405
$wsInd.Range("B5").Value = 937023.7502758912
$wsInd.Range("B6").Value = 1012244.4912122109
$wsInd.Range("B7").Value = 1029987.4107759548
$wsInd.Range("B8").Value = 1066009.9648907091
$wsInd.Range("B9").Value = 1117690.1755524143
$wsInd.Range("B10").Value = 1122695.6884051363
$wsInd.Range("B11").Value = 1164527.82235003
$wsInd.Range("B12").Value = 1185286.8803847947
$wsInd.Range("B13").Value = 1185384.6609280694
$wsInd.Range("B14").Value = 1142475.8593310704
$wsInd.Range("B15").Value = 1138414.8840670069
$wsInd.Range("B16").Value = 1150661.2880954784
$wsInd.Range("B17").Value = 1180336.1719917036
$wsInd.Range("B18").Value = 1084249.3612125
$wsInd.Range("B19").Value = 1114894.2239582008
$wsInd.Range("B20").Value = 1132674.5153484996
$wsInd.Range("B21").Value = 1135277.4777463179
$wsInd.Range("B22").Value = 1151487.9784192315
$wsInd.Range("B23").Value = 1139603.894875345
$wsInd.Range("B24").Value = 1131106.8244867884
$wsInd.Range("B25").Value = 1128459.7750324016
$wsInd.Range("B26").Value = 1124606.6865500417
$wsInd.Range("B27").Value = 1127104.6963498918
$wsInd.Range("B28").Value = 1127520.6918648488
$wsInd.Range("B29").Value = 1125516.2250042984
$wsInd.Range("B30").Value = 1060659.6030455015
$wsInd.Range("B31").Value = 1061055.2880089344
$wsInd.Range("B32").Value = 1065499.9892827943
$wsInd.Range("B33").Value = 1072620.7577635364
$wsInd.Range("B34").Value = 1135770.117537596
$wsInd.Range("B35").Value = 1134113.0043248208
$wsInd.Range("B36").Value = 1121631.3875612176
$wsInd.Range("B37").Value = 1099177.6432509285
$wsInd.Range("B38").Value = 1205597.8928211923
$wsInd.Range("B39").Value = 1177701.6191501068
$wsInd.Range("B40").Value = 1158954.1335180823
$wsInd.Range("B41").Value = 1151427.737185688
$wsInd.Range("B42").Value = 1186814.5152755342
$wsInd.Range("B43").Value = 1193422.6170181904
$wsInd.Range("B44").Value = 1200653.74433066
$wsInd.Range("B45").Value = 1206272.5518997707
$wsInd.Range("B46").Value = 1228429.9151612595
$wsInd.Range("B47").Value = 1230550.1387321956
$wsInd.Range("B48").Value = 1232325.0163338482
$wsInd.Range("B49").Value = 1234690.7683964064
$wsInd.Range("B50").Value = 1300344.1974023965
$wsInd.Range("B51").Value = 1306926.6112458303
$wsInd.Range("B52").Value = 1316261.5987988391
$wsInd.Range("B53").Value = 1327880.6201230814
$wsInd.Range("B54").Value = 1337795.8748758314
$wsInd.Range("B55").Value = 1349740.085446263
$wsInd.Range("B56").Value = 1359463.9926887949
$wsInd.Range("B57").Value = 1365918.71086328
$wsInd.Range("B58").Value = 1474849.0195391856
$wsInd.Range("B59").Value = 1473138.895730603
$wsInd.Range("B60").Value = 1467989.4191516854
$wsInd.Range("B61").Value = 1461678.0430899644
$wsInd.Range("B62").Value = 1417724.470969155
$wsInd.Range("B63").Value = 1420877.1067407792
$wsInd.Range("B64").Value = 1432915.5349530273
$wsInd.Range("B65").Value = 1452984.7709540613
$wsInd.Range("B66").Value = 1644827.0404468474
$wsInd.Range("B67").Value = 1669898.0884452055
$wsInd.Range("B68").Value = 1688890.592094805
$wsInd.Range("B69").Value = 1701119.6272407076
$wsInd.Range("B70").Value = 1737779.2511036643
$wsInd.Range("B71").Value = 1747416.3848487742
$wsInd.Range("B72").Value = 1761462.435642438
$wsInd.Range("B73").Value = 1779879.0653348616
$wsInd.Range("B74").Value = 1736877.0942185703
$wsInd.Range("B75").Value = 1961089.3789478329
$wsInd.Range("B76").Value = 1896613.4740806296
$wsInd.Range("B77").Value = 1882097.3194965683
$wsInd.Range("B78").Value = 2121171.282744103
$wsInd.Range("B79").Value = 2495512.1187999076
$wsInd.Range("B80").Value = 2455316.920402798
$wsInd.Range("B81").Value = 2447788.924153872
$wsInd.Range("B82").Value = 3301722.7553905826
$wsInd.Range("B83").Value = 3566202.869179457
$wsInd.Range("B84").Value = 3438202.2537474823
$wsInd.Range("B85").Value = 3550332.6303950218
$wsInd.Range("B86").Value = 4195514.327750179
$wsInd.Range("B87").Value = 4615164.88253978
$wsInd.Range("B88").Value = 4354898.611406246
$wsInd.Range("B89").Value = 4567212.970366247
$wsInd.Range("B90").Value = 4566786.9895859035
$wsInd.Range("B91").Value = 4963362.975722045
$wsInd.Range("B92").Value = 4657293.994309506
$wsInd.Range("B93").Value = 4942615.222752693
$wsInd.Range("B94").Value = 4712152.733703494
$wsInd.Range("B95").Value = 5322462.211726913
$wsInd.Range("B96").Value = 4960914.712555511
$wsInd.Range("B97").Value = 5283660.326120422
$wsInd.Range("B98").Value = 4868591.750611411
$wsInd.Range("B99").Value = 5221195.610913908
$wsInd.Range("B100").Value = 5332671.56935372
$wsInd.Range("B101").Value = 5241119.043827821

# --- Sheet "conso" (sheet4): update column B values ---
$wsConso = $wb.Worksheets.Item("conso")
$wsConso.Range("B2").Value = 295816.9965489449
$wsConso.Range("B3").Value = 370692.36013302696
$wsConso.Range("B4").Value = 411280.59703711804
$wsConso.Range("B5").Value = 392216.0132769617
$wsConso.Range("B6").Value = 363139.4666677071
$wsConso.Range("B7").Value = 372227.31880729255
$wsConso.Range("B8").Value = 398809.0646645876
$wsConso.Range("B9").Value = 363031.42136684904
$wsConso.Range("B10").Value = 367111.8240339916
$wsConso.Range("B11").Value = 400407.17527251417
$wsConso.Range("B12").Value = 411706.4597179249
$wsConso.Range("B13").Value = 427712.69536105386
$wsConso.Range("B14").Value = 456188.29991258855
$wsConso.Range("B15").Value = 472760.6916413203
$wsConso.Range("B16").Value = 517838.6058062659
$wsConso.Range("B17").Value = 505303.32377236773
$wsConso.Range("B18").Value = 604630.7490842085
$wsConso.Range("B19").Value = 644390.5454066086
$wsConso.Range("B20").Value = 694716.7074850892
$wsConso.Range("B21").Value = 734777.0
$wsConso.Range("B22").Value = 848205.0
$wsConso.Range("B23").Value = 885602.6366256276
$wsConso.Range("B24").Value = 961786.9927690151
$wsConso.Range("B25").Value = 1002140.7827081879
$wsConso.Range("B26").Value = 980210.5432055574

# --- Sheet "VA" (sheet3) values recalc automatically via formula pro!B - conso!B ---

# --- Update selections / view state to match the edited ranges ---
$wsPro.Range("B102:B113").Select()

$wsInd.Application.Goto($wsInd.Range("A96"), $False)
$wsInd.Range("B102:B113").Select()

$wsVA = $wb.Worksheets.Item("VA")
$wsVA.Range("B102:B113").Select()

$wsConso.Range("B102:B113").Select()

$wsPro.Select()
$wsPro.Range("B102:B113").Select()
